# Add a new "2022-Q4" sheet (fund holdings detail) right before the
# existing "2022-Q3" sheet, and insert a matching summary row at the top
# of the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet before "2022-Q3", cloning the header
#    / index-column look from the sheet that is about to become Q3 so
#    the new sheet matches the existing quarterly-report styling.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3HeaderFormat = $q3.Range("B1:H1")
$q3IndexFormat = $q3.Range("A2")

$newSheet = $wb.Worksheets.Add($q3)
$newSheet.Name = "2022-Q4"

# Clone header formatting (bold / border / centred) onto row 1, then
# overwrite with the real header text.
$q3HeaderFormat.Copy()
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Data rows: index, code, name, scale, position, ratio, value, rank.
# Numeric-looking strings are prefixed with a literal apostrophe so the
# interop layer stores them as text (matches the source export, which
# keeps fund codes / percentages as strings) instead of silently
# re-parsing them as numbers.
$data = @(
    @(0, "'002685", "中欧丰泓沪港深灵活配置混合A", "'54.86", "'93.29", "'3.90", "'2.1395", 9),
    @(1, "'002686", "中欧丰泓沪港深灵活配置混合C", "'16.87", "'93.29", "'3.90", "'0.6579", 9),
    @(2, "'005504", "汇添富沪港深大盘价值混合A", "'4.33", "'92.73", "'4.74", "'0.2052", 5),
    @(3, "'006205", "汇添富沪港深优势精选定期开放混合", "'0.55", "'91.51", "'4.60", "'0.0253", 7),
    @(4, "'517880", "华泰柏瑞中证沪港深品牌消费50ETF", "'0.48", "'92.59", "'5.11", "'0.0245", 6),
    @(5, "'015118", "汇添富沪港深大盘价值混合C", "'0.02", "'92.73", "'4.74", "'0.0009", 5),
    @(6, "'015119", "汇添富沪港深大盘价值混合D", "'0.00", "'92.73", "'4.74", 0, 5)
)

$rowNum = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($rowNum, 2).Value = $row[1]
    $newSheet.Cells.Item($rowNum, 3).Value = $row[2]
    $newSheet.Cells.Item($rowNum, 4).Value = $row[3]
    $newSheet.Cells.Item($rowNum, 5).Value = $row[4]
    $newSheet.Cells.Item($rowNum, 6).Value = $row[5]
    $newSheet.Cells.Item($rowNum, 7).Value = $row[6]
    $newSheet.Cells.Item($rowNum, 8).Value = $row[7]
    $rowNum++
}

# Clone the index-column format onto A2:A8, then fill in the 0-based
# row counter.
$q3IndexFormat.Copy()
$newSheet.Range("A2:A8").PasteSpecial($xlPasteFormats)
for ($r = 2; $r -le 8; $r++) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
}

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2) Insert a new top data row in "总计" summarising 2022-Q4.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$totalIndexFormat = $total.Range("A2")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$totalIndexFormat.Copy()
$total.Range("A2").PasteSpecial($xlPasteFormats)

$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 3.05

# The "A" column is a plain 0-based row counter; renumber it for every
# data row now that a new row sits on top.
for ($r = 2; $r -le 10; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

$total.Range("A1").Select()
